$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: Vianney worked 6 more hours (row 9, column B)
$ws.Range("B9").Value = 6

# Scroll the view down so row 12 becomes the top visible row, and move the
# active selection to B10 (matches the saved view state after the edit).
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B10").Select()

$excel.Calculate()
